# Daily scrape update - 2025-11-07 03:21:40 UTC
# Replaces the opportunity listing rows with the newly scraped data and
# drops the now-unused "premium" (yellow) highlight formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; A="1329388"; C="Sales Intern Home Appliances"; D="Panamá, Provincia de Panamá, Panamá"; E="No"; F="2 applicants"; G="6 - 18 Months"; H="Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"},
    @{Row=3; A="1329370"; C="Sales Deduction & Budget Control Intern"; D="Panamá, Provincia de Panamá, Panamá"; E="No"; F="2 applicants"; G="6 - 18 Months"; H="Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"},
    @{Row=4; A="1329350"; C="Technical Production Supervisor"; D="Accra, Ghana"; E="No"; F="1 applicant"; G="6 - 18 Months"; H="Dolcevita Chocolate LTD"},
    @{Row=5; A="1329349"; C="Copywriter"; D="Yerevan, Armenia"; E="No"; F="9 applicants"; G="6 - 18 Months"; H="TCF Armenia"},
    @{Row=6; A="1329337"; C="Junior Software Engineer"; D="Yerevan, Armenia"; E="No"; F="21 applicants"; G="6 - 18 Months"; H="Information Systems Agency of Armenia (ISAA)"},
    @{Row=7; A="1329336"; C="Junior SOC Analyst"; D="Yerevan, Armenia"; E="No"; F="10 applicants"; G="6 - 18 Months"; H="Information Systems Agency of Armenia (ISAA)"},
    @{Row=8; A="1329333"; C="CSIRT Analyst –  Incident Response (IR)"; D="Yerevan, Armenia"; E="No"; F="8 applicants"; G="6 - 18 Months"; H="Information Systems Agency of Armenia (ISAA)"},
    @{Row=9; A="1329257"; C="SALES MANAGER"; D="Denizli, Kumkısık, Denizli, Türkiye"; E="No"; F="5 applicants"; G="6 - 18 Months"; H="ALKA METAL"},
    @{Row=10; A="1329177"; C="Social Media Intern – Intercultural Engagement Program"; D="Porto Alegre, RS, Brazil"; E="No"; F="8 applicants"; G="9 - 12 Weeks"; H="ESCOLA GIORDANO BRUNO LTDA"},
    @{Row=11; A="1329121"; C="Sales and Distribution"; D="Accra, Ghana"; E="No"; F="1 applicant"; G="6 - 18 Months"; H="Dolcevita Chocolate LTD"},
    @{Row=12; A="1328962"; C="Social Media Manager"; D="Hong Kong"; E="No"; F="8 applicants"; G="6 - 18 Months"; H="Wong's Limited"},
    @{Row=13; A="1328930"; C="Business Development"; D="Istanbul, İstanbul, Türkiye"; E="No"; F="4 applicants"; G="3 - 6 Months"; H="Dentekay Dental Clinic"},
    @{Row=14; A="1328566"; C="HR Intern"; D="Santiago, Región Metropolitana, Chile"; E="No"; F="83 applicants"; G="6 - 18 Months"; H="Boehringer Ingelheim in Chile"},
    @{Row=15; A="1328558"; C="Flutter Developer"; D="Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"; E="No"; F="5 applicants"; G="9 - 12 Weeks"; H="Techno square"},
    @{Row=16; A="1328274"; C="Web Designer"; D="Tunis, Tunisia"; E="No"; F="19 applicants"; G="9 - 12 Weeks"; H="La fabrique"},
    @{Row=17; A="1328272"; C="Graphic Designer"; D="Tunis, Tunisia"; E="No"; F="10 applicants"; G="9 - 12 Weeks"; H="La fabrique"},
    @{Row=18; A="1328271"; C="Video Editor"; D="Tunis, Tunisia"; E="No"; F="4 applicants"; G="9 - 12 Weeks"; H="La fabrique"},
    @{Row=19; A="1327553"; C="Digital illustrator"; D="Cairo, Cairo Governorate, Egypt"; E="No"; F="1 applicant"; G="3 - 6 Months"; H="Perfect design firm"},
    @{Row=20; A="1327242"; C="Sales Manager"; D="Cairo, Cairo Governorate, Egypt"; E="No"; F="11 applicants"; G="9 - 12 Weeks"; H="MZ creatives"},
    @{Row=21; A="1325922"; C="Taste Hungary| Associate Technical Support - Polish"; D="Budapest, Magyarország"; E="No"; F="30 applicants"; G="6 - 18 Months"; H="Tech Mahindra Kft."},
    @{Row=22; A="1325612"; C="Retail Management Trainee Intern"; D="Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"; E="No"; F="23 applicants"; G="6 - 18 Months"; H="Vimigo"},
    @{Row=23; A="1325556"; C="Retail Management Trainee Intern"; D="Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"; E="No"; F="27 applicants"; G="6 - 18 Months"; H="Big Bath Sdn Bhd"},
    @{Row=24; A="1325555"; C="Talent Acquisition Intern"; D="Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"; E="No"; F="38 applicants"; G="6 - 18 Months"; H="Big Bath Sdn Bhd"},
    @{Row=25; A="1321840"; C="Marketing&Sales"; D="Gaziantep, Türkiye"; E="No"; F="67 applicants"; G="6 - 18 Months"; H="Hak Makarna"},
    @{Row=26; A="1321833"; C="Digital Marketing"; D="Gaziantep, Türkiye"; E="No"; F="65 applicants"; G="6 - 18 Months"; H="MACRO SIGN REKLAM İÇ VE DIŞ TİCARET LİMİTED ŞİRKETİ"},
    @{Row=27; A="1321254"; C="Marketing&Sales"; D="Gaziantep, Türkiye"; E="No"; F="101 applicants"; G="6 - 18 Months"; H="Baharoğlu Gıda"},
    @{Row=28; A="1314934"; C="Social Media Marketing Executive"; D="Petaling Jaya, Selangor, Malaysia"; E="No"; F="141 applicants"; G="6 - 18 Months"; H="iWisers SDN BHD"},
    @{Row=29; A="1310418"; C="MARKETING"; D="Gaziantep, Türkiye"; E="No"; F="56 applicants"; G="6 - 18 Months"; H="HÜNER"}
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $r.A
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
}

# The three rows that used to be flagged "Yes" / highlighted premium
# (E2, E11, E12) are no longer premium opportunities, and the
# corresponding yellow-fill style is dropped entirely - reset any
# lingering cell formatting back to the workbook default.
$ws.Range("E2").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Style = "Normal"

# Widen a few columns to fit the refreshed copy.
$ws.Columns.Item(3).ColumnWidth = 56.18   # -> stored width 57
$ws.Columns.Item(4).ColumnWidth = 69.18   # -> stored width 70
$ws.Columns.Item(8).ColumnWidth = 59.18   # -> stored width 60
